$wb = $excel.ActiveWorkbook

# --- Sheet "Installation" (sheet1): remove the Electric_boiler column (column B) ---
$wsInstallation = $wb.Worksheets.Item("Installation")
$wsInstallation.Range("B1").EntireColumn.Delete()

# --- Sheet "Capacity" (sheet2): remove the Electric_boiler column (column B) ---
$wsCapacity = $wb.Worksheets.Item("Capacity")
$wsCapacity.Range("B1").EntireColumn.Delete()

# Update the remaining capacity values (price policy update) on sheet "Capacity"
$wsCapacity.Range("B2").Value = 22.150119221048680
$wsCapacity.Range("F2").Value = 2500
$wsCapacity.Range("B3").Value = 38.319706252414221
$wsCapacity.Range("C3").Value = 545.08100032685832

# --- Sheet "Storage_capacity" (sheet3): update Hot_water_tank capacity value ---
$wsStorage = $wb.Worksheets.Item("Storage_capacity")
$wsStorage.Range("B2").Value = 897.60282631709015

# Save so that the now-unused "Electric_boiler" shared string is dropped
# and every shared-string index across the workbook is renumbered.
$wb.Save()
